$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1099
$ws1.Range("F3").Value = 413
$ws1.Range("F4").Value = 1507
$ws1.Range("F5").Value = 8770
$ws1.Range("F6").Value = 95
$ws1.Range("F9").Value = 290
$ws1.Range("F10").Value = 157
$ws1.Range("F11").Value = 20
$ws1.Range("F12").Value = 8
$ws1.Range("F13").Value = 3616
$ws1.Range("F17").Value = 1722
$ws1.Range("F20").Value = 312
$ws1.Range("F22").Value = 2399
$ws1.Range("F23").Value = 69

# Sheet "全部类型" (sheet4.xml): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1099
$ws4.Range("F3").Value = 413
$ws4.Range("F4").Value = 1507
$ws4.Range("F5").Value = 8770
$ws4.Range("F6").Value = 95
$ws4.Range("F9").Value = 290
$ws4.Range("F10").Value = 157
$ws4.Range("F11").Value = 20
$ws4.Range("F12").Value = 8
$ws4.Range("F13").Value = 3616
$ws4.Range("F17").Value = 1723
$ws4.Range("F20").Value = 312
$ws4.Range("F22").Value = 2399
$ws4.Range("F24").Value = 69
